# Fixed issue with missing attachments:
# update the OneDrive attachments folder path (Config sheet, row "OneDriveFolder")
# and leave the selection on the next cell down, matching the authored workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "C:\Users\GM00061060\OneDrive - ICU Medical Inc\Blackline Reconciliations - 2023\"

$ws.Range("B5").Select()
